$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

$ws.Range("C13").Value = "assay"
$ws.Range("D13").Value = "DNA Sequencing"

$ws.Range("C14").Value = "http://purl.obolibrary.org/obo/OBI_0000070"
$ws.Range("D14").Value = "http://purl.obolibrary.org/obo/NCIT_C153598"

$ws.Range("D15").Value = "NCIT"
